$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-18
$i0values = @(8, 9, 9, 9, 7, 8, 9, 8, 8, 7, 8, 3, 7, 10, 7, 8, 8)
$ifvalues = @(8, 9, 9, 9, 8, 8, 9, 8, 8, 8, 8, 4, 8, 10, 7, 8, 8)

for ($i = 0; $i -lt 17; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $i0values[$i]
    $ws.Cells.Item($row, 10).Value = $ifvalues[$i]
}
